$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet, this becomes the "valid" data sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ValidDataTypes"

# H2 used to hold a static date literal; make it a live TODAY() formula.
$ws1.Range("H2").Formula = "=TODAY()"

# Selection on ValidDataTypes moved to N12.
$ws1.Range("N12").Select() | Out-Null

# --- Build the "invalid" data sheet as a copy of the valid one, placed right after it ---
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "InvalidDataTypes"

# Row 3: a text value "Nineteen Ninety Nine" sneaks into the numeric Year column,
# and the same text lands in the date ("Date Watched") column too.
$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "Nineteen Ninety Nine"
$ws2.Range("H3").Value = "Nineteen Ninety Nine"

# Row 5: B5 (Year) becomes an actual date value, and H5 (Date Watched) becomes text.
# Borrow the existing date format (from H4, a cell that keeps the plain date style)
# so this reuses the workbook's existing date style instead of minting a new one.
$ws2.Range("H4").Copy() | Out-Null
$ws2.Range("B5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws2.Range("B5").Value = 44330
$ws2.Range("H5").Value = "Makoto Shinkai"

# Row 7: a text value "$42 billion" sneaks into the numeric Year column,
# and the same text lands in the date column too.
$ws2.Range("B7").NumberFormat = "0.00"
$ws2.Range("B7").Value = "$42 billion"
$ws2.Range("H7").Value = "$42 billion"

# Rows 8-10 are blanked out entirely (still styled, just no content) to represent
# missing/invalid rows.
$ws2.Range("A8:H10").ClearContents() | Out-Null

# Selection on InvalidDataTypes moved to A22.
$ws2.Range("A22").Select() | Out-Null

# Restore ValidDataTypes as the active tab.
$ws1.Activate()
$ws1.Range("N12").Select() | Out-Null
